$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.256564333333333
$ws.Range("H2").Value = 15.769693
$ws.Range("I2").Value = 0.003747859920520347
$ws.Range("J2").Value = 0.003747859920520347
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 112.513392
$ws.Range("N2").Value = 337.540176
$ws.Range("O2").Value = 0.3275312977368564
$ws.Range("P2").Value = 0.3275312977368564
$ws.Range("Q2").Value = 591.433883409552
$ws.Range("R2").Value = 5322.904950685967
$ws.Range("S2").Value = 0.001227541423503981
$ws.Range("T2").Value = 0.001227541423503981

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.256564333333333
$ws.Range("H3").Value = 15.769693
$ws.Range("I3").Value = 0.003747859920520347
$ws.Range("J3").Value = 0.003747859920520347
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 106.314466
$ws.Range("N3").Value = 318.943398
$ws.Range("O3").Value = 0.3094859589441663
$ws.Range("P3").Value = 0.3094859589441664
$ws.Range("Q3").Value = 558.8488300929794
$ws.Range("R3").Value = 5029.639470836814
$ws.Range("S3").Value = 0.001159910021490647
$ws.Range("T3").Value = 0.001159910021490647

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.256564333333333
$ws.Range("H4").Value = 15.769693
$ws.Range("I4").Value = 0.003747859920520347
$ws.Range("J4").Value = 0.003747859920520347
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 124.6916553333333
$ws.Range("N4").Value = 374.074966
$ws.Range("O4").Value = 0.3629827433189773
$ws.Range("P4").Value = 0.3629827433189773
$ws.Range("Q4").Value = 655.4497080894931
$ws.Range("R4").Value = 5899.047372805438
$ws.Range("S4").Value = 0.00136040847552572
$ws.Range("T4").Value = 0.00136040847552572

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1312.703450666667
$ws.Range("H5").Value = 3938.110352
$ws.Range("I5").Value = 0.93593996730609
$ws.Range("J5").Value = 0.9359399673060897
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 112.513392
$ws.Range("N5").Value = 337.540176
$ws.Range("O5").Value = 0.3275312977368564
$ws.Range("P5").Value = 0.3275312977368564
$ws.Range("Q5").Value = 147696.7179246113
$ws.Range("R5").Value = 1329270.461321502
$ws.Range("S5").Value = 0.3065496320955546
$ws.Range("T5").Value = 0.3065496320955545

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1312.703450666667
$ws.Range("H6").Value = 3938.110352
$ws.Range("I6").Value = 0.93593996730609
$ws.Range("J6").Value = 0.9359399673060897
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 106.314466
$ws.Range("N6").Value = 318.943398
$ws.Range("O6").Value = 0.3094859589441663
$ws.Range("P6").Value = 0.3094859589441664
$ws.Range("Q6").Value = 139559.366373984
$ws.Range("R6").Value = 1256034.297365856
$ws.Range("S6").Value = 0.2896602782958969
$ws.Range("T6").Value = 0.2896602782958969

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1312.703450666667
$ws.Range("H7").Value = 3938.110352
$ws.Range("I7").Value = 0.93593996730609
$ws.Range("J7").Value = 0.9359399673060897
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 124.6916553333333
$ws.Range("N7").Value = 374.074966
$ws.Range("O7").Value = 0.3629827433189773
$ws.Range("P7").Value = 0.3629827433189773
$ws.Range("Q7").Value = 163683.1662254054
$ws.Range("R7").Value = 1473148.496028648
$ws.Range("S7").Value = 0.3397300569146384
$ws.Range("T7").Value = 0.3397300569146384

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 84.59089266666666
$ws.Range("H8").Value = 253.772678
$ws.Range("I8").Value = 0.06031217277338979
$ws.Range("J8").Value = 0.06031217277338978
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 112.513392
$ws.Range("N8").Value = 337.540176
$ws.Range("O8").Value = 0.3275312977368564
$ws.Range("P8").Value = 0.3275312977368564
$ws.Range("Q8").Value = 9517.608266234591
$ws.Range("R8").Value = 85658.47439611131
$ws.Range("S8").Value = 0.01975412421779785
$ws.Range("T8").Value = 0.01975412421779785

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 84.59089266666666
$ws.Range("H9").Value = 253.772678
$ws.Range("I9").Value = 0.06031217277338979
$ws.Range("J9").Value = 0.06031217277338978
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 106.314466
$ws.Range("N9").Value = 318.943398
$ws.Range("O9").Value = 0.3094859589441663
$ws.Range("P9").Value = 0.3094859589441664
$ws.Range("Q9").Value = 8993.235582319981
$ws.Range("R9").Value = 80939.12024087984
$ws.Range("S9").Value = 0.01866577062677878
$ws.Range("T9").Value = 0.01866577062677878

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 84.59089266666666
$ws.Range("H10").Value = 253.772678
$ws.Range("I10").Value = 0.06031217277338979
$ws.Range("J10").Value = 0.06031217277338978
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 124.6916553333333
$ws.Range("N10").Value = 374.074966
$ws.Range("O10").Value = 0.3629827433189773
$ws.Range("P10").Value = 0.3629827433189773
$ws.Range("Q10").Value = 10547.77843273099
$ws.Range("R10").Value = 94930.00589457895
$ws.Range("S10").Value = 0.02189227792881316
$ws.Range("T10").Value = 0.02189227792881316

